# Daily attendance processing - 2025-11-01 22:19:35
#
# Normalizes the "Recorded By" column (G) on the active sheet: each cell
# holds a comma-separated list of recorder names/emails (e.g.
# "System, dnasr281@gmail.com"). The attendance sync job re-appends any
# exact-case "System" entry to the end of the list instead of leaving it
# first, so the human recorders sort ahead of the automated one. Lowercase
# "system" tokens (a separate, distinct entry) are left where they are.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $text = $cell.Value2

    if ($text -eq $null) {
        continue
    }
    if ($text.GetType().Name -ne "String") {
        continue
    }
    if ($text.IndexOf(",") -lt 0) {
        continue
    }

    $parts = $text -split ", "

    # NOTE: PowerShell's "case-sensitive" operators (-ceq/-cne/-cmatch) are
    # NOT reliably case-sensitive in this host, so case-sensitive equality
    # is done via [string]::Equals(...) (ordinal, case-sensitive by default).
    $kept = @()
    $systemParts = @()
    foreach ($part in $parts) {
        if ($part.Equals("System")) {
            $systemParts += $part
        } else {
            $kept += $part
        }
    }

    if ($systemParts.Count -eq 0) {
        continue
    }

    $ordered = $kept + $systemParts
    $newText = [string]::Join(", ", $ordered)

    if (-not $newText.Equals($text)) {
        $cell.Value2 = $newText
    }
}
